$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.087389945983887
$ws.Range("B1").Value = 3.094913244247437
$ws.Range("C1").Value = 6.58001708984375
$ws.Range("D1").Value = 1.830418825149536
$ws.Range("E1").Value = 1.24502694606781
